$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2
$ws.Range("F2").Value = 31
$ws.Range("G2").Value = "adam"
$ws.Range("I2").Value = 64
$ws.Range("J2").Value = 28.54411922053105
$ws.Range("K2").Value = 1337.11523739707
$ws.Range("L2").Value = 36.56658635143661
$ws.Range("M2").Value = 0.168715769032278

# Row 3 (new)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "GRU"
$ws.Range("D3").Value = 30
$ws.Range("E3").Value = 60
$ws.Range("F3").Value = 31
$ws.Range("G3").Value = "<keras.src.optimizers.legacy.adam.Adam object at 0x7b66a52a5960>"
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 32
$ws.Range("J3").Value = 31.70121541473221
$ws.Range("K3").Value = 1506.28299943318
$ws.Range("L3").Value = 38.81086187439258
$ws.Range("M3").Value = 0.1902124552115729

# Row 4 (new)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "GRU"
$ws.Range("D4").Value = 40
$ws.Range("E4").Value = 60
$ws.Range("F4").Value = 31
$ws.Range("G4").Value = "adam"
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 64
$ws.Range("J4").Value = 47.20802092508893
$ws.Range("K4").Value = 2880.612931639214
$ws.Range("L4").Value = 53.6713418095655
$ws.Range("M4").Value = 0.285709775643143

# Row 5 (new)
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "GRU"
$ws.Range("D5").Value = 40
$ws.Range("E5").Value = 60
$ws.Range("F5").Value = 31
$ws.Range("G5").Value = "<keras.src.optimizers.legacy.adam.Adam object at 0x7b6627fe2dd0>"
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 32
$ws.Range("J5").Value = 28.32222051535054
$ws.Range("K5").Value = 1382.150259845443
$ws.Range("L5").Value = 37.17728150154934
$ws.Range("M5").Value = 0.1797625645883431

# Apply the same style as A2 (bold/border header-like style) to A3:A5
$ws.Range("A2").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)  # xlPasteFormats
